# SECURDE Documentation.docx edit:
# - Add a new row of content describing that the application now logs user
#   activity (fills the previously-empty row that only held the hidden
#   _GoBack bookmark).
# - Merge the two runs in the "hidden buttons" vulnerability cell into one.
# - Resize the table / columns / row heights to their new (slightly larger)
#   values, matching Word's relayout after the edits.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellParagraphXml($cell, [string]$innerParagraph) {
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerParagraph + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$cell.Range.InsertXML($xml)
}

# --- Fill in the new "Application does not log user activities" row (row 8) ---
# Cell 1 is set first so the pre-existing hidden _GoBack bookmark (which lived
# in that cell's otherwise-empty paragraph) is discarded.
Set-CellParagraphXml $t.Cell(8, 1) '<w:p><w:r><w:t>Application does not log user activities</w:t></w:r></w:p>'
Set-CellParagraphXml $t.Cell(8, 2) '<w:p><w:r><w:t>LogWrite.java, Main.java, Frame.java and SQLite.java</w:t></w:r></w:p>'
Set-CellParagraphXml $t.Cell(8, 3) '<w:p><w:r><w:t>Users could be accessing or modifying content not assigned to their role without anyone knowing</w:t></w:r></w:p>'
Set-CellParagraphXml $t.Cell(8, 4) '<w:p><w:r><w:t>User logins, logs out, registers or them accessing content panels gets logged into a .txt file.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

# --- Merge the two runs in row 6, column 3 into a single run ---
Set-CellParagraphXml $t.Cell(6, 3) '<w:p><w:r><w:t>Users could modify data or see sensitive information that would normally be hidden</w:t></w:r></w:p>'

# --- Resize columns (this also updates each tcW and the tblGrid) ---
$t.Columns.Item(1).Width = 162.8
$t.Columns.Item(2).Width = 162.8
$t.Columns.Item(3).Width = 162.85
$t.Columns.Item(4).Width = 162.85

# --- Resize the overall table (preferred) width ---
$t.PreferredWidth = 651.3

# --- Resize each row's height to match the new layout ---
$t.Rows.Item(1).Height = 16
$t.Rows.Item(2).Height = 35.85
$t.Rows.Item(3).Height = 37.9
$t.Rows.Item(4).Height = 35.85
$t.Rows.Item(5).Height = 37.9
$t.Rows.Item(6).Height = 35.85
$t.Rows.Item(7).Height = 37.9
$t.Rows.Item(8).Height = 37.9
$t.Rows.Item(9).Height = 35.85
$t.Rows.Item(10).Height = 37.9
$t.Rows.Item(11).Height = 35.85
